# Work log update for 张非凡 — add June 2019 vue/tomcat integration entries.
# Cells are written in the same order the original author entered them so the
# shared-string table builds up in the same sequence as the real workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Row 35: extend the existing 2019.6.14 entry with a time range and add its content.
$ws.Range("B35").Value = "2019.6.14  15：00-17：00"
$ws.Range("C35").Value = "学习vue的基本模式和mvvm"

# Row 36
$ws.Range("B36").Value = "2019.6.17  15：00-17：00"
$ws.Range("C36").Value = "构建vue界面的基本框架"

# Row 37 date, row 39 date, row 40 content, row 37 content, row 38 date,
# row 40 date, row 38 content, row 39 content (matches original authoring order).
$ws.Range("B37").Value = "2019.6.18  18：00-20：00"
$ws.Range("B39").Value = "2019.6.20  18：00-20：00"
$ws.Range("C40").Value = "完成对数据的CRUD操作"
$ws.Range("C37").Value = "学习前后端的数据交互"
$ws.Range("B38").Value = "2019.6.19  14：00-17：00"
$ws.Range("B40").Value = "2019.6.21  20：00-22：00"
$ws.Range("C38").Value = "实现数据交互，对tomcat和vue先后启动"
$ws.Range("C39").Value = "实现数据交互，对其进行完善"

# Match the author's final selection/view state.
$ws.Range("A19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("C39").Select() | Out-Null
